# Auto-generated edit script applying the diff changes to Sophia_Profits workbook
# Each block targets one (sheet, row) pair; cells are set, cleared, or newly created
# to match the target OOXML exactly.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 20
$ws.Range("I29").Value = 20
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 60
$ws.Range("L29").ClearContents()
$ws.Range("N29").Value = 0
$ws.Range("M29").Value = 221

$ws.Range("H38").Value = 1002.4
$ws.Range("I38").Value = 878
$ws.Range("J38").Value = 1500
$ws.Range("K38").Value = 2634
$ws.Range("L38").Value = 4500
$ws.Range("M38").Value = -2262
$ws.Range("N38").Value = -5244

$ws.Range("H58").Value = 4014.375
$ws.Range("I58").Value = 315
$ws.Range("J58").Value = 4542.857
$ws.Range("K58").Value = 945
$ws.Range("L58").Value = 13628.571
$ws.Range("M58").Value = -795
$ws.Range("N58").Value = -13928.571

$ws.Range("H70").Value = 1221.875
$ws.Range("I70").Value = 887.5
$ws.Range("J70").Value = 1333.3334
$ws.Range("K70").Value = 2662.5
$ws.Range("L70").Value = 4000.0002
$ws.Range("M70").Value = -2392.5
$ws.Range("N70").Value = -4540.0002

$ws.Range("H73").Value = 1221.875
$ws.Range("I73").Value = 887.5
$ws.Range("J73").Value = 1333.3334
$ws.Range("K73").Value = 2662.5
$ws.Range("L73").Value = 4000.0002
$ws.Range("M73").Value = -1726.5
$ws.Range("N73").Value = -5872.0002

$ws.Range("H87").Value = 99354
$ws.Range("J87").Value = 99354
$ws.Range("L87").Value = 99354
$ws.Range("N87").Value = -101850

$ws.Range("H90").Value = 99354
$ws.Range("J90").Value = 99354
$ws.Range("L90").Value = 298062
$ws.Range("N90").Value = -310542

$ws.Range("H100").Value = 3999.25
$ws.Range("I100").Value = 4999
$ws.Range("K100").Value = 4999
$ws.Range("M100").Value = -4458

$ws.Range("H112").Value = 1818.5333
$ws.Range("I112").Value = 1549.4
$ws.Range("J112").Value = 1953.1
$ws.Range("K112").Value = 4648.200000000001
$ws.Range("L112").Value = 5859.299999999999
$ws.Range("M112").Value = -3540.200000000001
$ws.Range("N112").Value = -8075.299999999999

$ws.Range("H138").Value = 12864.108
$ws.Range("J138").Value = 13264.441
$ws.Range("L138").Value = 39793.323
$ws.Range("N138").Value = -50073.323

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1520.0714
$ws.Range("I2").Value = 1580.2
$ws.Range("K2").Value = 1580.2
$ws.Range("M2").Value = -1467.2

$ws.Range("H24").Value = 20355
$ws.Range("J24").Value = 20355
$ws.Range("L24").Value = 20355
$ws.Range("N24").Value = -21103

$ws.Range("H32").Value = 5508.4165
$ws.Range("I32").Value = 5991
$ws.Range("K32").Value = 5991
$ws.Range("M32").Value = -5704

$ws.Range("H61").Value = 1526.25
$ws.Range("I61").Value = 1526.25
$ws.Range("K61").Value = 1526.25
$ws.Range("M61").Value = -1314.25

$ws.Range("H97").Value = 697.375
$ws.Range("I97").Value = 697.375
$ws.Range("K97").Value = 697.375
$ws.Range("M97").Value = -201.375

$ws.Range("H100").Value = 20355
$ws.Range("J100").Value = 20355
$ws.Range("L100").Value = 20355
$ws.Range("N100").Value = -22519

$ws.Range("H102").Value = 2996.3333
$ws.Range("J102").Value = 3000
$ws.Range("L102").Value = 3000
$ws.Range("N102").Value = -6244

$ws.Range("H116").Value = 1520.0714
$ws.Range("I116").Value = 1580.2
$ws.Range("K116").Value = 1580.2
$ws.Range("M116").Value = 713.8

$ws.Range("H122").Value = 1604
$ws.Range("J122").Value = 1750
$ws.Range("L122").Value = 5250
$ws.Range("N122").Value = -10150

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").ClearContents()
$ws.Range("N130").Value = 0

$ws.Range("H136").Value = 1526.25
$ws.Range("I136").Value = 1526.25
$ws.Range("K136").Value = 4578.75
$ws.Range("M136").Value = -2028.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1520.0714
$ws.Range("I3").Value = 1580.2
$ws.Range("K3").Value = 1580.2
$ws.Range("M3").Value = -1466.2

$ws.Range("H63").Value = 89999
$ws.Range("J63").Value = 89999
$ws.Range("L63").Value = 89999
$ws.Range("N63").Value = -91371

$ws.Range("H66").Value = 89999
$ws.Range("J66").Value = 89999
$ws.Range("L66").Value = 269997
$ws.Range("N66").Value = -276861

$ws.Range("H94").Value = 5903.533
$ws.Range("I94").Value = 3379.4167
$ws.Range("J94").Value = 16000
$ws.Range("K94").Value = 3379.4167
$ws.Range("L94").Value = 16000
$ws.Range("M94").Value = -2928.4167
$ws.Range("N94").Value = -16902

$ws.Range("H107").Value = 1699.5714
$ws.Range("I107").Value = 1579.4
$ws.Range("K107").Value = 1579.4
$ws.Range("M107").Value = 340.5999999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 28747.25
$ws.Range("J86").Value = 4994.5
$ws.Range("L86").Value = 4994.5
$ws.Range("N86").Value = -7240.5

$ws.Range("H89").Value = 28747.25
$ws.Range("J89").Value = 4994.5
$ws.Range("L89").Value = 24972.5
$ws.Range("N89").Value = -36204.5

$ws.Range("H122").Value = 1383.5555
$ws.Range("I122").Value = 1383.5555
$ws.Range("K122").Value = 4150.666499999999
$ws.Range("M122").Value = -1700.666499999999

$ws.Range("H124").Value = 79663
$ws.Range("J124").Value = 79663
$ws.Range("L124").Value = 79663
$ws.Range("N124").Value = -84573

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1012
$ws.Range("I5").Value = 550.4286
$ws.Range("K5").Value = 1651.2858
$ws.Range("M5").Value = -1539.2858

$ws.Range("H23").Value = 184.14285
$ws.Range("J23").Value = 515
$ws.Range("L23").Value = 1545
$ws.Range("N23").Value = -2015

$ws.Range("H135").Value = 1012
$ws.Range("I135").Value = 550.4286
$ws.Range("K135").Value = 4953.8574
$ws.Range("M135").Value = -2418.8574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 326.77777
$ws.Range("I2").Value = 189.6
$ws.Range("J2").Value = 498.25
$ws.Range("K2").Value = 189.6
$ws.Range("L2").Value = 498.25
$ws.Range("M2").Value = -76.59999999999999
$ws.Range("N2").Value = -724.25

$ws.Range("H122").Value = 9466.666999999999
$ws.Range("I122").Value = 1700
$ws.Range("J122").Value = 25000
$ws.Range("K122").Value = 5100
$ws.Range("L122").Value = 75000
$ws.Range("M122").Value = -2650
$ws.Range("N122").Value = -79900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1872.75
$ws.Range("J16").Value = 1800
$ws.Range("L16").Value = 1800
$ws.Range("N16").Value = -2140

$ws.Range("H61").Value = 7847415
$ws.Range("I61").Value = 6376262
$ws.Range("K61").Value = 6376262
$ws.Range("M61").Value = -6376060

$ws.Range("H68").Value = 1629.3334
$ws.Range("J68").Value = 2000
$ws.Range("L68").Value = 2000
$ws.Range("N68").Value = -3498

$ws.Range("H71").Value = 1629.3334
$ws.Range("J71").Value = 2000
$ws.Range("L71").Value = 10000
$ws.Range("N71").Value = -17488

$ws.Range("H113").Value = 7847415
$ws.Range("I113").Value = 6376262
$ws.Range("K113").Value = 6376262
$ws.Range("M113").Value = -6374092

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3382.5454
$ws.Range("I122").Value = 2601.1428
$ws.Range("J122").Value = 4750
$ws.Range("K122").Value = 7803.428400000001
$ws.Range("L122").Value = 14250
$ws.Range("M122").Value = -5353.428400000001
$ws.Range("N122").Value = -19150
